# no-op for now
$p = $ppt.ActivePresentation
